$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header labels in row 1 (A1:L1) but keep their formatting/style.
$ws.Range("A1:L1").ClearContents()

# Remove the data in row 4 entirely (cells + the row itself collapse out of
# sheetData, subsequent rows keep their original row numbers).
$ws.Rows.Item(4).ClearContents()

# Update the selection to span the (now empty) header row.
$ws.Range("A1:L1").Select()
